$wb = $excel.ActiveWorkbook

$wsMenu = $wb.Worksheets.Item("Menu List")
$wsMod  = $wb.Worksheets.Item("Modifier List")

# --- Menu List sheet: rows 22-24 edits ---

# A22:A24 - replace placeholder random string, drop the border around it
$wsMenu.Range("A22:A24").Value = "pFiMb9HUdJ"
$wsMenu.Range("A22:A24").Borders.LineStyle = -4142

# B22:D24 - drop the border (text/values stay the same)
$wsMenu.Range("B22:D24").Borders.LineStyle = -4142

# F22:G22 - strip the red/centered/bordered style back to plain default
$wsMenu.Range("F22:G22").Style = "Normal"

# F23:I23 - values become "3" (as text) and lose their style
$wsMenu.Range("F23:I23").Value = "'3"
$wsMenu.Range("F23:I23").Style = "Normal"

# F24:K24 - values stay the same, just lose their style
$wsMenu.Range("F24:K24").Style = "Normal"

# --- Active sheet / selection bookkeeping ---
# Modifier List keeps its old selection (H35) but is no longer the active tab;
# select it first so that activating Menu List afterwards is what "sticks".
$wsMod.Range("H35").Select()

# Menu List becomes the active tab, with G33 selected
$wsMenu.Activate()
$wsMenu.Range("G33").Select()
